$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture formatting of the old footer rows (31-32) before we clear them ---
$ws.Range("B31:C32").Copy()
$ws.Range("B50:C51").PasteSpecial(-4122)
$ws.Range("H31:J32").Copy()
$ws.Range("H50:J51").PasteSpecial(-4122)
$ws.Range("B50:C50").Merge()
$ws.Range("B51:C51").Merge()
$ws.Range("H50:J50").Merge()
$ws.Range("H51:J51").Merge()

# --- Capture the special highlighted-row formatting (old row 26) for the new last data row (45) ---
$ws.Range("B26:J26").Copy()
$ws.Range("B45:J45").PasteSpecial(-4122)

# --- Remove the old footer rows 31-32 (unmerge first, then clear) ---
$ws.Range("B31:C31").UnMerge()
$ws.Range("B32:C32").UnMerge()
$ws.Range("H31:J31").UnMerge()
$ws.Range("H32:J32").UnMerge()
$ws.Range("B31:J32").Clear()

# --- Apply the normal data-row formatting (old row 16) across all new data rows 17-44 ---
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J44").PasteSpecial(-4122)

# --- Update the three summary header values ---
$ws.Range("E11").Value = 1965065
$ws.Range("C13").Value = 13
$ws.Range("F13").Value = 8

# --- Fill in the 30 detail rows (16-45) ---
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1143331660"
$ws.Cells.Item(16, 4).Value = "CARMEN JULIANA BELTRAN BEDOYA"
$ws.Cells.Item(16, 5).Value = "2504"
$ws.Cells.Item(16, 6).Value = 100000
$ws.Cells.Item(16, 7).Value = 2500000

$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1143331660"
$ws.Cells.Item(17, 4).Value = "CARMEN JULIANA BELTRAN BEDOYA"
$ws.Cells.Item(17, 5).Value = "2504"
$ws.Cells.Item(17, 6).Value = 100000
$ws.Cells.Item(17, 7).Value = 2500000

$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1143350848"
$ws.Cells.Item(18, 4).Value = "LILIBETH LICONA CHIQUILLO"
$ws.Cells.Item(18, 5).Value = "2507"
$ws.Cells.Item(18, 6).Value = 98622
$ws.Cells.Item(18, 7).Value = 2465540

$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "1143350848"
$ws.Cells.Item(19, 4).Value = "LILIBETH LICONA CHIQUILLO"
$ws.Cells.Item(19, 5).Value = "2504"
$ws.Cells.Item(19, 6).Value = 98622
$ws.Cells.Item(19, 7).Value = 2465540

$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "1014205331"
$ws.Cells.Item(20, 4).Value = "MILEIDYS PAJARO GARCES"
$ws.Cells.Item(20, 5).Value = "2507"
$ws.Cells.Item(20, 6).Value = 56940
$ws.Cells.Item(20, 7).Value = 1423500

$ws.Cells.Item(21, 2).Value = "CC"
$ws.Cells.Item(21, 3).Value = "1014205331"
$ws.Cells.Item(21, 4).Value = "MILEIDYS PAJARO GARCES"
$ws.Cells.Item(21, 5).Value = "2504"
$ws.Cells.Item(21, 6).Value = 56940
$ws.Cells.Item(21, 7).Value = 1423500

$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "1143340357"
$ws.Cells.Item(22, 4).Value = "WILIAN OVIEDO MENDOZA"
$ws.Cells.Item(22, 5).Value = "2502"
$ws.Cells.Item(22, 6).Value = 13867
$ws.Cells.Item(22, 7).Value = 1423500

$ws.Cells.Item(23, 2).Value = "CC"
$ws.Cells.Item(23, 3).Value = "1143369813"
$ws.Cells.Item(23, 4).Value = "LIZ KARIME ROMERO PALOMINO"
$ws.Cells.Item(23, 5).Value = "2507"
$ws.Cells.Item(23, 6).Value = 91600
$ws.Cells.Item(23, 7).Value = 2290000

$ws.Cells.Item(24, 2).Value = "CC"
$ws.Cells.Item(24, 3).Value = "1143369813"
$ws.Cells.Item(24, 4).Value = "LIZ KARIME ROMERO PALOMINO"
$ws.Cells.Item(24, 5).Value = "2504"
$ws.Cells.Item(24, 6).Value = 91600
$ws.Cells.Item(24, 7).Value = 2290000

$ws.Cells.Item(25, 2).Value = "CC"
$ws.Cells.Item(25, 3).Value = "1128059061"
$ws.Cells.Item(25, 4).Value = "CINDY PAOLA PARRA PERNA"
$ws.Cells.Item(25, 5).Value = "2507"
$ws.Cells.Item(25, 6).Value = 56940
$ws.Cells.Item(25, 7).Value = 1423500

$ws.Cells.Item(26, 2).Value = "CC"
$ws.Cells.Item(26, 3).Value = "1128059061"
$ws.Cells.Item(26, 4).Value = "CINDY PAOLA PARRA PERNA"
$ws.Cells.Item(26, 5).Value = "2504"
$ws.Cells.Item(26, 6).Value = 56940
$ws.Cells.Item(26, 7).Value = 1423500

$ws.Cells.Item(27, 2).Value = "CC"
$ws.Cells.Item(27, 3).Value = "1143352257"
$ws.Cells.Item(27, 4).Value = "YUBENIS HERNANDEZ RODRIGUEZ"
$ws.Cells.Item(27, 5).Value = "2507"
$ws.Cells.Item(27, 6).Value = 56940
$ws.Cells.Item(27, 7).Value = 1423500

$ws.Cells.Item(28, 2).Value = "CC"
$ws.Cells.Item(28, 3).Value = "1143352257"
$ws.Cells.Item(28, 4).Value = "YUBENIS HERNANDEZ RODRIGUEZ"
$ws.Cells.Item(28, 5).Value = "2504"
$ws.Cells.Item(28, 6).Value = 56940
$ws.Cells.Item(28, 7).Value = 1423500

$ws.Cells.Item(29, 2).Value = "CC"
$ws.Cells.Item(29, 3).Value = "1043640354"
$ws.Cells.Item(29, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(29, 5).Value = "2507"
$ws.Cells.Item(29, 6).Value = 52000
$ws.Cells.Item(29, 7).Value = 1300000

$ws.Cells.Item(30, 2).Value = "CC"
$ws.Cells.Item(30, 3).Value = "1043640354"
$ws.Cells.Item(30, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(30, 5).Value = "2506"
$ws.Cells.Item(30, 6).Value = 52000
$ws.Cells.Item(30, 7).Value = 1300000

$ws.Cells.Item(31, 2).Value = "CC"
$ws.Cells.Item(31, 3).Value = "1043640354"
$ws.Cells.Item(31, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(31, 5).Value = "2505"
$ws.Cells.Item(31, 6).Value = 52000
$ws.Cells.Item(31, 7).Value = 1300000

$ws.Cells.Item(32, 2).Value = "CC"
$ws.Cells.Item(32, 3).Value = "1043640354"
$ws.Cells.Item(32, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(32, 5).Value = "2504"
$ws.Cells.Item(32, 6).Value = 52000
$ws.Cells.Item(32, 7).Value = 1300000

$ws.Cells.Item(33, 2).Value = "CC"
$ws.Cells.Item(33, 3).Value = "1043640354"
$ws.Cells.Item(33, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(33, 5).Value = "2503"
$ws.Cells.Item(33, 6).Value = 52000
$ws.Cells.Item(33, 7).Value = 1300000

$ws.Cells.Item(34, 2).Value = "CC"
$ws.Cells.Item(34, 3).Value = "1043640354"
$ws.Cells.Item(34, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(34, 5).Value = "2502"
$ws.Cells.Item(34, 6).Value = 52000
$ws.Cells.Item(34, 7).Value = 1300000

$ws.Cells.Item(35, 2).Value = "CC"
$ws.Cells.Item(35, 3).Value = "1043640354"
$ws.Cells.Item(35, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(35, 5).Value = "2501"
$ws.Cells.Item(35, 6).Value = 52000
$ws.Cells.Item(35, 7).Value = 1300000

$ws.Cells.Item(36, 2).Value = "CC"
$ws.Cells.Item(36, 3).Value = "1043640354"
$ws.Cells.Item(36, 4).Value = "DANIELA SOFIA MENDOZA CHAVEZ"
$ws.Cells.Item(36, 5).Value = "2412"
$ws.Cells.Item(36, 6).Value = 52000
$ws.Cells.Item(36, 7).Value = 1300000

$ws.Cells.Item(37, 2).Value = "CC"
$ws.Cells.Item(37, 3).Value = "1143384006"
$ws.Cells.Item(37, 4).Value = "AGUSTIN WILMAN NAAR PEREZ"
$ws.Cells.Item(37, 5).Value = "2507"
$ws.Cells.Item(37, 6).Value = 76207
$ws.Cells.Item(37, 7).Value = 1905190

$ws.Cells.Item(38, 2).Value = "CC"
$ws.Cells.Item(38, 3).Value = "1143384006"
$ws.Cells.Item(38, 4).Value = "AGUSTIN WILMAN NAAR PEREZ"
$ws.Cells.Item(38, 5).Value = "2504"
$ws.Cells.Item(38, 6).Value = 76207
$ws.Cells.Item(38, 7).Value = 1905190

$ws.Cells.Item(39, 2).Value = "CC"
$ws.Cells.Item(39, 3).Value = "1049931644"
$ws.Cells.Item(39, 4).Value = "KATIA PAOLA PINTO PEREZ"
$ws.Cells.Item(39, 5).Value = "2507"
$ws.Cells.Item(39, 6).Value = 88000
$ws.Cells.Item(39, 7).Value = 2200000

$ws.Cells.Item(40, 2).Value = "CC"
$ws.Cells.Item(40, 3).Value = "1049931644"
$ws.Cells.Item(40, 4).Value = "KATIA PAOLA PINTO PEREZ"
$ws.Cells.Item(40, 5).Value = "2504"
$ws.Cells.Item(40, 6).Value = 88000
$ws.Cells.Item(40, 7).Value = 2200000

$ws.Cells.Item(41, 2).Value = "CC"
$ws.Cells.Item(41, 3).Value = "1032376335"
$ws.Cells.Item(41, 4).Value = "ANDRES RICARDO OSORIO MALPICA"
$ws.Cells.Item(41, 5).Value = "2507"
$ws.Cells.Item(41, 6).Value = 56940
$ws.Cells.Item(41, 7).Value = 1423500

$ws.Cells.Item(42, 2).Value = "CC"
$ws.Cells.Item(42, 3).Value = "1032376335"
$ws.Cells.Item(42, 4).Value = "ANDRES RICARDO OSORIO MALPICA"
$ws.Cells.Item(42, 5).Value = "2504"
$ws.Cells.Item(42, 6).Value = 56940
$ws.Cells.Item(42, 7).Value = 1423500

$ws.Cells.Item(43, 2).Value = "CC"
$ws.Cells.Item(43, 3).Value = "1193032566"
$ws.Cells.Item(43, 4).Value = "ENILDA ISABEL MARRUGO ELLES"
$ws.Cells.Item(43, 5).Value = "2507"
$ws.Cells.Item(43, 6).Value = 56940
$ws.Cells.Item(43, 7).Value = 1423500

$ws.Cells.Item(44, 2).Value = "PPT"
$ws.Cells.Item(44, 3).Value = "5488453"
$ws.Cells.Item(44, 4).Value = "NETHYARI ALEJANDRA GUZMAN HERRERA"
$ws.Cells.Item(44, 5).Value = "2507"
$ws.Cells.Item(44, 6).Value = 56940
$ws.Cells.Item(44, 7).Value = 1423500

$ws.Cells.Item(45, 2).Value = "PPT"
$ws.Cells.Item(45, 3).Value = "5488453"
$ws.Cells.Item(45, 4).Value = "NETHYARI ALEJANDRA GUZMAN HERRERA"
$ws.Cells.Item(45, 5).Value = "2504"
$ws.Cells.Item(45, 6).Value = 56940
$ws.Cells.Item(45, 7).Value = 1423500

# --- New footer rows 50-51 ---
$ws.Cells.Item(50, 2).Value = "___________________________________"
$ws.Cells.Item(50, 8).Value = "___________________________________"
$ws.Cells.Item(51, 2).Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Cells.Item(51, 8).Value = "FIRMA DEL REPRESENTANTE LEGAL"

Write-Host "done"
